$d = $word.ActiveDocument

# Original paragraph text is "Version 1." laid out as:
#   [proofErr spellStart] "Version" [proofErr spellEnd] " 1." [bookmarkStart _GoBack][bookmarkEnd]
#
# Target layout is "Version 2." laid out as:
#   [proofErr spellStart] "Versi" "on" [proofErr spellEnd] " 2" [bookmarkStart _GoBack][bookmarkEnd] "."
#
# i.e. the word "Version" gets split into two runs ("Versi" + "on"), the
# " 1." run becomes " 2", and the final "." is moved into its own run,
# placed after the (pre-existing) _GoBack bookmark.
#
# Range.InsertXML is used (instead of simply assigning .Text, or toggling a
# character property to force a run split) because it lets us replace the
# targeted range with an exact set of <w:r> elements without leaving any
# run-formatting (w:rPr) residue behind.

# Step 1: split "Version" (chars 0-6) into "Versi" (0-4) + "on" (5-6) by
# rewriting just the "on" tail (chars 5,7) in place. This stays inside the
# existing run/proofErr span, so the spellStart/spellEnd markers are left
# exactly where they were -- just now wrapping two runs instead of one.
$xmlOn = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rSplitVersion = $d.Range(5, 7)
$rSplitVersion.InsertXML($xmlOn)

# Step 2: turn " 1" (chars 7,9) into " 2", leaving the trailing "."
# (now immediately following, still before the bookmark) untouched for now.
$xmlTwo = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> 2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rVersionNumber = $d.Range(7, 9)
$rVersionNumber.InsertXML($xmlTwo)

# Step 3: the "." now sits right before the _GoBack bookmark (chars 7,8).
# Remove it from there...
$rOldDot = $d.Range(7, 8)
$rOldDot.Delete()

# Step 4: ...and re-insert it as its own run after the bookmark, at the
# (now) end of the paragraph.
$rEnd = $d.Range(9, 9)
$rEnd.InsertAfter(".")
